$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 874.13336
$ws.Cells.Item(33, 9).Value = 656.0909
$ws.Cells.Item(33, 10).Value = 1473.75
$ws.Cells.Item(33, 11).Value = 656.0909
$ws.Cells.Item(33, 12).Value = 1473.75
$ws.Cells.Item(33, 13).Value = -427.0909
$ws.Cells.Item(33, 14).Value = -1931.75

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 1249.0714
$ws.Cells.Item(107, 9).Value = 1338.909
$ws.Cells.Item(107, 11).Value = 1338.909
$ws.Cells.Item(107, 13).Value = 581.0909999999999

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 18114.715
$ws.Cells.Item(116, 9).Value = 18114.715
$ws.Cells.Item(116, 11).Value = 18114.715
$ws.Cells.Item(116, 13).Value = -14672.715

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1400.8
$ws.Cells.Item(137, 9).Value = 1037.625
$ws.Cells.Item(137, 11).Value = 3112.875
$ws.Cells.Item(137, 13).Value = -562.875

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3092.4255
$ws.Cells.Item(138, 9).Value = 2556.45
$ws.Cells.Item(138, 10).Value = 3489.4443
$ws.Cells.Item(138, 11).Value = 7669.349999999999
$ws.Cells.Item(138, 12).Value = 10468.3329
$ws.Cells.Item(138, 13).Value = -2529.349999999999
$ws.Cells.Item(138, 14).Value = -20748.3329

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4482.232
$ws.Cells.Item(32, 9).Value = 2898.9783
$ws.Cells.Item(32, 10).Value = 11765.2
$ws.Cells.Item(32, 11).Value = 2898.9783
$ws.Cells.Item(32, 12).Value = 11765.2
$ws.Cells.Item(32, 13).Value = -2611.9783
$ws.Cells.Item(32, 14).Value = -12339.2

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 254223.19
$ws.Cells.Item(61, 9).Value = 254223.19
$ws.Cells.Item(61, 11).Value = 254223.19
$ws.Cells.Item(61, 13).Value = -254011.19

# ARM row 121
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 9).Value = 0
$ws.Cells.Item(121, 11).Value = 0
$ws.Cells.Item(121, 13).ClearContents()

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 4954.676
$ws.Cells.Item(132, 9).Value = 5281.7095
$ws.Cells.Item(132, 10).Value = 3265
$ws.Cells.Item(132, 11).Value = 15845.1285
$ws.Cells.Item(132, 12).Value = 9795
$ws.Cells.Item(132, 13).Value = -13315.1285
$ws.Cells.Item(132, 14).Value = -14855

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 254223.19
$ws.Cells.Item(136, 9).Value = 254223.19
$ws.Cells.Item(136, 11).Value = 762669.5700000001
$ws.Cells.Item(136, 13).Value = -760119.5700000001

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1633.875
$ws.Cells.Item(99, 9).Value = 865.7143
$ws.Cells.Item(99, 11).Value = 865.7143
$ws.Cells.Item(99, 13).Value = 632.2857

# BSM row 112
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(112, 8).Value = 49000
$ws.Cells.Item(112, 10).Value = 49000
$ws.Cells.Item(112, 12).Value = 49000
$ws.Cells.Item(112, 14).Value = -51954

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2350.56
$ws.Cells.Item(134, 9).Value = 1723.0889
$ws.Cells.Item(134, 10).Value = 7997.8
$ws.Cells.Item(134, 11).Value = 5169.2667
$ws.Cells.Item(134, 12).Value = 23993.4
$ws.Cells.Item(134, 13).Value = -2634.2667
$ws.Cells.Item(134, 14).Value = -29063.4

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 20889.8
$ws.Cells.Item(22, 9).Value = 34333
$ws.Cells.Item(22, 10).Value = 725
$ws.Cells.Item(22, 11).Value = 34333
$ws.Cells.Item(22, 12).Value = 725
$ws.Cells.Item(22, 13).Value = -33983
$ws.Cells.Item(22, 14).Value = -1425

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 5611.75
$ws.Cells.Item(86, 9).Value = 4916
$ws.Cells.Item(86, 11).Value = 4916
$ws.Cells.Item(86, 13).Value = -3793

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 5611.75
$ws.Cells.Item(89, 9).Value = 4916
$ws.Cells.Item(89, 11).Value = 24580
$ws.Cells.Item(89, 13).Value = -18964

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2377.1667
$ws.Cells.Item(99, 9).Value = 2368
$ws.Cells.Item(99, 10).Value = 2404.6667
$ws.Cells.Item(99, 11).Value = 2368
$ws.Cells.Item(99, 12).Value = 2404.6667
$ws.Cells.Item(99, 13).Value = -870
$ws.Cells.Item(99, 14).Value = -5400.6667

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 2377.1667
$ws.Cells.Item(126, 9).Value = 2368
$ws.Cells.Item(126, 10).Value = 2404.6667
$ws.Cells.Item(126, 11).Value = 7104
$ws.Cells.Item(126, 12).Value = 7214.000100000001
$ws.Cells.Item(126, 13).Value = -4634
$ws.Cells.Item(126, 14).Value = -12154.0001

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 884.5333000000001
$ws.Cells.Item(132, 9).Value = 855.9167
$ws.Cells.Item(132, 11).Value = 2567.7501
$ws.Cells.Item(132, 13).Value = -37.7501000000002

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 626.8
$ws.Cells.Item(92, 9).Value = 626.8
$ws.Cells.Item(92, 11).Value = 1880.4
$ws.Cells.Item(92, 13).Value = -632.3999999999999

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 951.6923
$ws.Cells.Item(122, 9).Value = 581.2
$ws.Cells.Item(122, 10).Value = 1456.909
$ws.Cells.Item(122, 11).Value = 5230.8
$ws.Cells.Item(122, 12).Value = 13112.181
$ws.Cells.Item(122, 13).Value = -2780.8
$ws.Cells.Item(122, 14).Value = -18012.181

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 1346.75
$ws.Cells.Item(140, 9).Value = 1103.2
$ws.Cells.Item(140, 11).Value = 3309.6
$ws.Cells.Item(140, 13).Value = 1870.4

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1390.2
$ws.Cells.Item(102, 9).Value = 1399.6316
$ws.Cells.Item(102, 11).Value = 1399.6316
$ws.Cells.Item(102, 13).Value = 222.3684000000001

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 803.0625
$ws.Cells.Item(107, 9).Value = 486.42307
$ws.Cells.Item(107, 10).Value = 2175.1667
$ws.Cells.Item(107, 11).Value = 486.42307
$ws.Cells.Item(107, 12).Value = 2175.1667
$ws.Cells.Item(107, 13).Value = 1433.57693
$ws.Cells.Item(107, 14).Value = -6015.1667

# GSM row 109
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3599.5789
$ws.Cells.Item(113, 9).Value = 2578.5833
$ws.Cells.Item(113, 11).Value = 2578.5833
$ws.Cells.Item(113, 13).Value = -408.5832999999998

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 5000.3706
$ws.Cells.Item(126, 9).Value = 5000.3706
$ws.Cells.Item(126, 11).Value = 15001.1118
$ws.Cells.Item(126, 13).Value = -12531.1118

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3154.1667
$ws.Cells.Item(16, 9).Value = 919.44446
$ws.Cells.Item(16, 10).Value = 5388.8887
$ws.Cells.Item(16, 11).Value = 919.44446
$ws.Cells.Item(16, 12).Value = 5388.8887
$ws.Cells.Item(16, 13).Value = -749.44446
$ws.Cells.Item(16, 14).Value = -5728.8887

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 181.69048
$ws.Cells.Item(55, 9).Value = 124.37037
$ws.Cells.Item(55, 10).Value = 284.86667
$ws.Cells.Item(55, 11).Value = 124.37037
$ws.Cells.Item(55, 12).Value = 284.86667
$ws.Cells.Item(55, 13).Value = 48.62963000000001
$ws.Cells.Item(55, 14).Value = -630.86667

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1223.4706
$ws.Cells.Item(100, 9).Value = 1172.5
$ws.Cells.Item(100, 11).Value = 2345
$ws.Cells.Item(100, 13).Value = -1804

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2070.9333
$ws.Cells.Item(126, 9).Value = 2167.4
$ws.Cells.Item(126, 10).Value = 1878
$ws.Cells.Item(126, 11).Value = 6502.200000000001
$ws.Cells.Item(126, 12).Value = 5634
$ws.Cells.Item(126, 13).Value = -4032.200000000001
$ws.Cells.Item(126, 14).Value = -10574
